# Updates cryptocurrency price (column D) and 1h volume change (column E) values.
# Values must remain plain text (matching the source sheet's inline-string cells),
# so we briefly force a Text number format while assigning, then restore "General"
# to avoid leaving the cell visually marked as text-formatted.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "326.04" }
    @{ Cell = "E2"; Value = "-0.37%" }
    @{ Cell = "D3"; Value = "44.33" }
    @{ Cell = "E3"; Value = "-1.82%" }
    @{ Cell = "D4"; Value = "5.500" }
    @{ Cell = "E4"; Value = "-1.80%" }
    @{ Cell = "D5"; Value = "0.07999" }
    @{ Cell = "E5"; Value = "-1.35%" }
    @{ Cell = "D6"; Value = "2.026" }
    @{ Cell = "E6"; Value = "6.72%" }
    @{ Cell = "D7"; Value = "4.298" }
    @{ Cell = "E7"; Value = "-1.09%" }
    @{ Cell = "E8"; Value = "-6.86%" }
    @{ Cell = "D10"; Value = "0.1143" }
    @{ Cell = "E10"; Value = "-3.16%" }
    @{ Cell = "D11"; Value = "0.1837" }
    @{ Cell = "E11"; Value = "-3.49%" }
    @{ Cell = "D12"; Value = "12.16" }
    @{ Cell = "E12"; Value = "41.39%" }
    @{ Cell = "D13"; Value = "0.09751" }
    @{ Cell = "E13"; Value = "-3.84%" }
    @{ Cell = "D14"; Value = "0.04553" }
    @{ Cell = "E14"; Value = "9.50%" }
    @{ Cell = "D16"; Value = "0.001264" }
    @{ Cell = "E16"; Value = "-0.86%" }
    @{ Cell = "D17"; Value = "0.04081" }
    @{ Cell = "E17"; Value = "-4.66%" }
    @{ Cell = "D18"; Value = "0.005876" }
    @{ Cell = "E18"; Value = "-0.72%" }
    @{ Cell = "D19"; Value = "3.365" }
    @{ Cell = "E19"; Value = "-6.39%" }
    @{ Cell = "E20"; Value = "-0.20%" }
    @{ Cell = "D21"; Value = "0.1405" }
    @{ Cell = "E21"; Value = "2.07%" }
    @{ Cell = "D23"; Value = "0.001242" }
    @{ Cell = "E23"; Value = "-0.02%" }
    @{ Cell = "D24"; Value = "0.004305" }
    @{ Cell = "E24"; Value = "-5.65%" }
    @{ Cell = "D25"; Value = "0.0001189" }
    @{ Cell = "E25"; Value = "-3.84%" }
    @{ Cell = "D26"; Value = "0.0003739" }
    @{ Cell = "E26"; Value = "-6.70%" }
    @{ Cell = "D38"; Value = "0.02546" }
    @{ Cell = "E38"; Value = "-5.07%" }
    @{ Cell = "D39"; Value = "0.05529" }
    @{ Cell = "E39"; Value = "-1.84%" }
    @{ Cell = "D40"; Value = "0.007530" }
    @{ Cell = "E40"; Value = "-2.35%" }
    @{ Cell = "E41"; Value = "-0.46%" }
    @{ Cell = "D42"; Value = "0.007588" }
    @{ Cell = "E42"; Value = "-33.05%" }
    @{ Cell = "D43"; Value = "0.002013" }
    @{ Cell = "E43"; Value = "-2.49%" }
    @{ Cell = "D44"; Value = "0.008384" }
    @{ Cell = "E44"; Value = "-3.60%" }
    @{ Cell = "D45"; Value = "0.00007094" }
    @{ Cell = "E45"; Value = "-0.21%" }
    @{ Cell = "D46"; Value = "0.00000000749" }
    @{ Cell = "E46"; Value = "-0.64%" }
    @{ Cell = "E47"; Value = "0.93%" }
    @{ Cell = "D48"; Value = "0.004219" }
    @{ Cell = "E48"; Value = "21.79%" }
    @{ Cell = "D49"; Value = "0.00002097" }
    @{ Cell = "E49"; Value = "-0.64%" }
    @{ Cell = "D50"; Value = "0.0001998" }
    @{ Cell = "E50"; Value = "-0.64%" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.NumberFormat = "General"
}

Write-Output ("Applied {0} cell updates" -f $updates.Count)
